# Add a new week of "Apio" (Mercado Mayorista Lo Valledor de Santiago) data.
# Two new rows are inserted right after the existing row 330, pushing all the
# subsequent rows (old 331..364) down by two (new 333..366).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 331/332; everything below shifts down.
$ws.Rows("331:332").Insert()

# --- New row 331 (Primera) ---
$ws.Cells.Item(331, 1).Value = 6
$ws.Cells.Item(331, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(331, 3).Value = "Metropolitana"
$ws.Cells.Item(331, 4).Value = 44449
$ws.Cells.Item(331, 5).Value = 13
$ws.Cells.Item(331, 6).Value = 100112017
$ws.Cells.Item(331, 7).Value = "Apio"
$ws.Cells.Item(331, 8).Value = "Americana (o)"
$ws.Cells.Item(331, 9).Value = "Primera"
$ws.Cells.Item(331, 10).Value = 2100
$ws.Cells.Item(331, 11).Value = 7000
$ws.Cells.Item(331, 12).Value = 8000
$ws.Cells.Item(331, 13).Value = 7571
$ws.Cells.Item(331, 14).Value = '$/docena de matas'
$ws.Cells.Item(331, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(331, 16).Value = 1262
$ws.Cells.Item(331, 17).Value = 6
$ws.Cells.Item(331, 18).Value = "Hortaliza"

# --- New row 332 (Segunda) ---
$ws.Cells.Item(332, 1).Value = 6
$ws.Cells.Item(332, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(332, 3).Value = "Metropolitana"
$ws.Cells.Item(332, 4).Value = 44449
$ws.Cells.Item(332, 5).Value = 13
$ws.Cells.Item(332, 6).Value = 100112017
$ws.Cells.Item(332, 7).Value = "Apio"
$ws.Cells.Item(332, 8).Value = "Americana (o)"
$ws.Cells.Item(332, 9).Value = "Segunda"
$ws.Cells.Item(332, 10).Value = 800
$ws.Cells.Item(332, 11).Value = 6000
$ws.Cells.Item(332, 12).Value = 6000
$ws.Cells.Item(332, 13).Value = 6000
$ws.Cells.Item(332, 14).Value = '$/docena de matas'
$ws.Cells.Item(332, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(332, 16).Value = 1000
$ws.Cells.Item(332, 17).Value = 6
$ws.Cells.Item(332, 18).Value = "Hortaliza"
